# ---------------------------------------------------------------------------
# "Edited Chapter 2 of Literature Review"
#
# 1) The auto-updating "last saved" date placeholder (datetimeFigureOut
#    field) that lives on the slide master, every slide layout and the
#    notes master gets refreshed from 16/01/2023 -> 1/02/2023 (this is what
#    PowerPoint does by itself whenever the deck is opened/saved on a
#    different day).
# 2) The second picture on slide 13 ("Content Placeholder 8", the ALP-mass
#    plot next to "Picture 4") is resized/repositioned.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# ppPlaceholderDate
$ppPlaceholderDate = 16
$newDate = "1/02/2023"

function Update-DatePlaceholders {
    param($Container, $NewDate)

    $cnt = $Container.Shapes.Count
    for ($i = 1; $i -le $cnt; $i++) {
        $sh = $Container.Shapes.Item($i)
        $phType = -1
        try {
            $phType = $sh.PlaceholderFormat.Type
        } catch {
            $phType = -1
        }
        if ($phType -eq 16) {
            $sh.TextFrame.TextRange.Text = $NewDate
        }
    }
}

# Slide master.
Update-DatePlaceholders $p.SlideMaster $newDate

# Every custom (slide) layout hanging off the master.
$layoutCount = $p.SlideMaster.CustomLayouts.Count
for ($li = 1; $li -le $layoutCount; $li++) {
    Update-DatePlaceholders ($p.SlideMaster.CustomLayouts.Item($li)) $newDate
}

# Notes master.
Update-DatePlaceholders $p.NotesMaster $newDate

# ---------------------------------------------------------------------------
# Resize/reposition the second picture on slide 13.
# ---------------------------------------------------------------------------

function Set-ShapeEmu {
    param($Shape, $Prop, $TargetEmu)

    $emuPerPt = 12700.0
    $pts = $TargetEmu / $emuPerPt

    switch ($Prop) {
        "Left"   { $Shape.Left   = $pts }
        "Top"    { $Shape.Top    = $pts }
        "Width"  { $Shape.Width  = $pts }
        "Height" { $Shape.Height = $pts }
    }

    # Shape.Left/Top/Width/Height round-trip through single-precision
    # points, so the EMU value written back out can land 1 EMU away from
    # the target. Nudge it back in until it lands exactly.
    $actual = 0
    switch ($Prop) {
        "Left"   { $actual = $Shape.Left }
        "Top"    { $actual = $Shape.Top }
        "Width"  { $actual = $Shape.Width }
        "Height" { $actual = $Shape.Height }
    }
    $resultEmu = [math]::Round($actual * $emuPerPt)

    $tries = 0
    while ($resultEmu -ne $TargetEmu -and $tries -lt 30) {
        $diffEmu = $TargetEmu - $resultEmu
        $pts = $pts + ($diffEmu / $emuPerPt)
        switch ($Prop) {
            "Left"   { $Shape.Left   = $pts }
            "Top"    { $Shape.Top    = $pts }
            "Width"  { $Shape.Width  = $pts }
            "Height" { $Shape.Height = $pts }
        }
        $actual = 0
        switch ($Prop) {
            "Left"   { $actual = $Shape.Left }
            "Top"    { $actual = $Shape.Top }
            "Width"  { $actual = $Shape.Width }
            "Height" { $actual = $Shape.Height }
        }
        $resultEmu = [math]::Round($actual * $emuPerPt)
        $tries = $tries + 1
    }
}

$slide13 = $p.Slides.Item(13)

$pic = $null
for ($i = 1; $i -le $slide13.Shapes.Count; $i++) {
    $candidate = $slide13.Shapes.Item($i)
    if ($candidate.Name -eq "Content Placeholder 8") {
        $pic = $candidate
    }
}
if ($pic -eq $null) {
    $pic = $slide13.Shapes.Item(4)
}

Set-ShapeEmu $pic "Left"   5940583
Set-ShapeEmu $pic "Top"    2575956
Set-ShapeEmu $pic "Width"  6030469
Set-ShapeEmu $pic "Height" 2988253
